# Rename the first worksheet from "env" to "Environmental_variables".
$wb = $excel.ActiveWorkbook
$wsEnv = $wb.Worksheets.Item("env")
$wsEnv.Name = "Environmental_variables"

# Make the renamed sheet the active (tab-selected) sheet again and move its
# selection to I26. Activating it / selecting a range on it automatically
# clears the "tabSelected" flag on the other sheet and updates the
# workbook's active tab, matching the rest of the diff.
$wsEnv.Activate()
$null = $wsEnv.Range("I26").Select()
